$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A258").Value = "船"
$ws.Range("B258").Value = "ふね"
$ws.Range("C258").Value = "[1]"
$ws.Range("D258").Value = "船"

$ws.Range("A259").Value = "港"
$ws.Range("B259").Value = "みなと"
$ws.Range("C259").Value = "[0]"
$ws.Range("D259").Value = "港口、碼頭"

$ws.Range("A260").Value = "線路"
$ws.Range("B260").Value = "せんろ"
$ws.Range("C260").Value = "[1]"
$ws.Range("D260").Value = "鐵軌、鐵路"

$ws.Range("A261").Value = "汽車"
$ws.Range("B261").Value = "きしゃ"
$ws.Range("C261").Value = "[1][2]"
$ws.Range("D261").Value = "蒸汽火車"

$ws.Range("A262").Value = "列車"
$ws.Range("B262").Value = "れっしゃ"
$ws.Range("C262").Value = "[0][1]"
$ws.Range("D262").Value = "列車"

$ws.Range("A263").Value = "特急"
$ws.Range("B263").Value = "とっきゅう"
$ws.Range("C263").Value = "[0]"
$ws.Range("D263").Value = "特快車"

$ws.Range("A264").Value = "乗り物"
$ws.Range("B264").Value = "のりもの"
$ws.Range("C264").Value = "[0]"
$ws.Range("D264").Value = "交通工具"

$ws.Range("A265").Value = "新幹線"
$ws.Range("B265").Value = "しんかんせん"
$ws.Range("C265").Value = "[3]"
$ws.Range("D265").Value = "新幹線"

$ws.Range("A266").Value = "エスカレーター"
$ws.Range("B266").Value = "エスカレーター"
$ws.Range("C266").Value = "[4]"
$ws.Range("D266").Value = "電扶梯"

$ws.Range("A267").Value = "オートバイ"
$ws.Range("B267").Value = "オートバイ"
$ws.Range("C267").Value = "[3]"
$ws.Range("D267").Value = "摩托車"

$ws.Range("A268").Value = "機械"
$ws.Range("B268").Value = "きかい"
$ws.Range("C268").Value = "[2][1]"
$ws.Range("D268").Value = "機械"

$ws.Range("A269").Value = "電灯"
$ws.Range("B269").Value = "でんとう"
$ws.Range("C269").Value = "[0]"
$ws.Range("D269").Value = "電燈"

$ws.Range("A270").Value = "電球"
$ws.Range("B270").Value = "でんきゅう"
$ws.Range("C270").Value = "[0]"
$ws.Range("D270").Value = "電燈泡"

$ws.Range("A271").Value = "除湿機"
$ws.Range("B271").Value = "じょしつき"
$ws.Range("C271").Value = "[3]"
$ws.Range("D271").Value = "除濕機"

$ws.Range("A272").Value = "電子辞書"
$ws.Range("B272").Value = "でんしじしょ"
$ws.Range("C272").Value = "[4]"
$ws.Range("D272").Value = "電子辭典"

$ws.Range("A273").Value = "ベル"
$ws.Range("B273").Value = "ベル"
$ws.Range("C273").Value = "[1]"
$ws.Range("D273").Value = "電鈴"

$ws.Range("A274").Value = "ソフト"
$ws.Range("B274").Value = "ソフト"
$ws.Range("C274").Value = "[1]"
$ws.Range("D274").Value = "軟體"

$ws.Range("A275").Value = "ハード"
$ws.Range("B275").Value = "ハード"
$ws.Range("C275").Value = "[1]"
$ws.Range("D275").Value = "硬體"

$ws.Range("A276").Value = "ステレオ"
$ws.Range("B276").Value = "ステレオ"
$ws.Range("C276").Value = "[0]"
$ws.Range("D276").Value = "立體音響"

$ws.Range("A277").Value = "プリンター"
$ws.Range("B277").Value = "プリンター"
$ws.Range("C277").Value = "[0]"
$ws.Range("D277").Value = "印表機"

$ws.Range("A278").Value = "キーボード"
$ws.Range("B278").Value = "キーボード"
$ws.Range("C278").Value = "[3]"
$ws.Range("D278").Value = "鍵盤"

$ws.Range("A279").Value = "コンピューター"
$ws.Range("B279").Value = "コンピューター"
$ws.Range("C279").Value = "[3]"
$ws.Range("D279").Value = "電腦"

$ws.Range("A280").Value = "デジタルカメラ"
$ws.Range("B280").Value = "デジタルカメラ"
$ws.Range("C280").Value = "[5]"
$ws.Range("D280").Value = "數位相機"

[void]$ws.Range("D270").Select()
